$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "TestValue"
$v2 = $ws.Range("A1").Value2
Write-Host "V2:" $v2
